# Regenerate save_data to use K (strikeouts) instead of Strike# for the
# "young_danny" 2024 sheet. Column G ("K") values are recalculated from the
# underlying per-game data and rewritten below, row by row (rows 2-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 2
    21 = 2
    22 = 2
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 5
    28 = 2
    29 = 0
    30 = 2
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    43 = 2
    44 = 2
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
